$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "56.812.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "2.328.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'519.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'132.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.534"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "  -1.84%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'5.26"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.339"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'23.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "2.743.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "56.786.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "2.315.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.54%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'10.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'332.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.97%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'4.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "  +3.28%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'61.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.167"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.30%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'8.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.58%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'1.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.06%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'168.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "0.0₃0725"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.87%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'6.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'18.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.995"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.48%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'3.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.08%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.891"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.60"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.90%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'38.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.17%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'148.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.34%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "Bittensor"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'291.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.54%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.375"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "Filecoin"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'3.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'5.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.85%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0932"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0500"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.07%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.559"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'18.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.23%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0217"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'17.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.46%  "
$ws.Range("E50").Style = "Normal"
